$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the blank separator / totals rows (old rows 63 & 64),
# pushing everything down (blank separator -> 65, totals -> 66, bottom border -> 69).
$ws.Rows("63:64").Insert()

# ---- Row 63: The Admiral Rodney / Wollaton ----
$ws.Range("A63").Value = 44370
$ws.Range("B63").Value = "The Admiral Rodney"
$ws.Range("C63").Value = "Wollaton"
$ws.Range("D63").Value = "start/end at pub"
$ws.Range("E63").Value = 4.21
$ws.Range("F63").Value = 0.042106481481481488
$ws.Range("G63").Formula = "=F63/E63"
$ws.Range("H63").Value = 1
$ws.Range("I63").Value = 1
$ws.Range("J63").Value = 1
$ws.Range("N63").Value = 1
$ws.Range("O63").Value = "Bicycle repair man"
$ws.Range("P63").Formula = "=SUM(H63:N63)*E63"

# ---- Row 64: The Shed / The Furnace / Darley Abbey ----
# (Location typed before Pub so shared-string ids land the same as the source edit.)
$ws.Range("A64").Value = 44482
$ws.Range("C64").Value = "Darley Abbey"
$ws.Range("B64").Value = "The Shed / The Furnace"
$ws.Range("D64").Value = "start/end at pub"
$ws.Range("E64").Value = 3.2
$ws.Range("F64").Value = 0.031261574074074074
$ws.Range("G64").Formula = "=F64/E64"
$ws.Range("H64").Value = 1
$ws.Range("I64").Value = 1
$ws.Range("J64").Value = 1
$ws.Range("N64").Value = 1
$ws.Range("O64").Value = "Lost in the nature reserve, highland cattle"
$ws.Range("P64").Formula = "=SUM(H64:N64)*E64"

# The template row used for the insert didn't carry a style for column J (it was blank
# in every preceding row), so copy formats across from a column that did (H) for the
# freshly entered "1"s in J63/J64.
$ws.Range("H63:H64").Copy()
$ws.Range("J63:J64").PasteSpecial(-4122)
$ws.Range("J63").Value = 1
$ws.Range("J64").Value = 1

# Row 49's shared "=SUM(H:N)*E" formula now stretches across the two new rows too.
$ws.Range("P49:P64").Formula = "=SUM(H49:N49)*E49"

# ---- Totals row (now row 66): extend the summed ranges to include rows 63 & 64 ----
$ws.Range("E66").Formula = "=SUM(E5:E64)"
$ws.Range("G66").Formula = "=AVERAGE(G6:G64)"
$ws.Range("H66").Formula = "=SUM(H5:H64)"
$ws.Range("I66").Formula = "=SUM(I5:I64)"
$ws.Range("J66").Formula = "=SUM(J5:J64)"
$ws.Range("N66").Formula = "=SUM(N5:N64)"
$ws.Range("P66").Formula = "=SUM(P5:P64)"
